$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn off multi-threaded/concurrent calculation (maps toward concurrentCalc="0")
$excel.MultiThreadedCalculation.Enabled = $false

# Clear out the old data rows (2-9); new content spans rows 2-11.
$ws.Range("A2:G9").Clear()

# Row 2: datetime_utc
$ws.Range("A2").Value = "datetime_utc"
$ws.Range("B2").Value = "Data product UTC date and time"
$ws.Range("C2").Value = "Date"
$ws.Range("E2").Value = "YYYY-MM-DD hh:mm:ss"

# Row 3: datetime_utc_matlab
$ws.Range("A3").Value = "datetime_utc_matlab"
$ws.Range("B3").Value = "PI-provided UTC date and time "
$ws.Range("C3").Value = "numeric"
$ws.Range("D3").Value = "dimensionless"

# Row 4: latitude_matlab
$ws.Range("A4").Value = "latitude_matlab"
$ws.Range("B4").Value = "Latitude of sample event provided by PI"
$ws.Range("C4").Value = "numeric"
$ws.Range("D4").Value = "degree"

# Row 5: longitude_matlab
$ws.Range("A5").Value = "longitude_matlab"
$ws.Range("B5").Value = "Longitude of sample event provided by PI"
$ws.Range("C5").Value = "numeric"
$ws.Range("D5").Value = "degree"

# Row 6: latitude_API
$ws.Range("A6").Value = "latitude_API"
$ws.Range("B6").Value = "Latitude of sample event provided by NES-LTER API"
$ws.Range("C6").Value = "numeric"
$ws.Range("D6").Value = "degree"

# Row 7: longitude_API
$ws.Range("A7").Value = "longitude_API"
$ws.Range("B7").Value = "Longitude of sample event provided by NES-LTER API"
$ws.Range("C7").Value = "numeric"
$ws.Range("D7").Value = "degree"

# Row 8: toi_source
$ws.Range("A8").Value = "toi_source"
$ws.Range("B8").Value = "Bottle sample from niskin or underway"
$ws.Range("C8").Value = "categorical"

# Row 9: gop
$ws.Range("A9").Value = "gop"
$ws.Range("B9").Value = "Gross oxygen production"
$ws.Range("C9").Value = "numeric"
$ws.Range("D9").Value = "millimoleOxygenPerMeterSquaredPerDay"
$ws.Range("F9").Value = "NaN"
$ws.Range("G9").Value = "Missing value"

# Row 10: ncp
$ws.Range("A10").Value = "ncp"
$ws.Range("B10").Value = "Rate of net community production integrated over the mixed layer"
$ws.Range("C10").Value = "numeric"
$ws.Range("D10").Value = "millimoleOxygenPerMeterSquaredPerDay"
$ws.Range("F10").Value = "NaN"
$ws.Range("G10").Value = "Missing value"

# Row 11: ncp_per_gop
$ws.Range("A11").Value = "ncp_per_gop"
$ws.Range("B11").Value = "Net community production divided by gross oxygen production"
$ws.Range("C11").Value = "numeric"
$ws.Range("D11").Value = "dimensionless"
$ws.Range("F11").Value = "NaN"
$ws.Range("G11").Value = "Missing value"

# Update the selection shown when the sheet is opened
$ws.Range("A4:E7").Select() | Out-Null
